# "script changes of order placeing"
# Populate row 6 of the LoginTestData sheet with a new test case
# ("login_as_counter_and_place_an_order") and adjust the sheet's
# selection / column width bookkeeping accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row data (row 6)
$ws.Range("B6").Value = "login_as_counter_and_place_an_order"
$ws.Range("C6").Value = "Y"
$ws.Range("D6").Value = "Kilpauk!"
$ws.Range("E6").Value = "password"
$ws.Range("F6").Value = 1

# Widen column B so the longer test-case name fits, and select F6 as
# the last-touched / active cell (also clears the old topLeftCell scroll
# position pinned at D1).
$ws.Columns.Item(2).ColumnWidth = 35.15
$ws.Range("F6").Select()
